$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("text input")
$ws.Range("B2").Value = "This is the placeholder text for the prelim schedule individualized for each team. Useful pieces of information to include in this section are: expected start time and length of time for lunch, information on tiebreakers, where to report back after lunch..."
$ws.Range("B2").WrapText = $true
